$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 10635
$ws.Range("C2").Value = 551304.5

$ws.Range("B3").Value = 115875.8
$ws.Range("C3").Value = 1470196.49

$ws.Range("B4").Value = 525147.84
$ws.Range("C4").Value = 2263817.8

$ws.Range("B5").Value = 1337379.79
$ws.Range("C5").Value = 3512734.74

$ws.Range("B6").Value = 1604031.07
$ws.Range("C6").Value = 3520972.84

$ws.Range("B7").Value = 1067503.69
$ws.Range("C7").Value = 3402571.31

$ws.Range("B8").Value = 2685974.57
$ws.Range("C8").Value = 3817231.7

$ws.Range("B9").Value = 1080234.3
$ws.Range("C9").Value = 1452793.9
